$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right at the end of the text of the
# first (and only) paragraph. New paragraphs are being appended after it, and
# the bookmark needs to end up at the very end of the document once we're
# done, so remove it now and re-create it in the right spot afterwards.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$dash = [char]0x2013

# Append the new paragraphs (two blank lines interleaved with three new text
# paragraphs) after the existing "NGPS Editor Framework" paragraph.
$r = $d.Content
$r.Collapse(0)
$newText = "`r`rEditor $dash store space ( all main interface apps use this global variable for common global address space )" + `
           "`r`rEditor.dock $dash the upper dock that holds the presentation title and the loaded apps ( plus some util buttons like save and load )" + `
           "`rEditor.dock.UI $dash the builder class for the dock interface"
$r.InsertAfter($newText)

# Work around an engine edge-case: adding a *collapsed* bookmark exactly at
# the end-of-story position (or at the end of a paragraph's text, right
# before its paragraph mark, when that happens to be the last paragraph in
# the document) gets mis-anchored back to the start of the paragraph.
# To avoid this we temporarily add an extra trailing paragraph (so the
# target position is no longer the absolute end of the document), insert a
# one-character placeholder at the target spot (making it a true "interior"
# position), add the bookmark there, remove the placeholder again, and
# finally delete the temporary trailing paragraph.
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$lastContentParagraph = $d.Paragraphs($d.Paragraphs.Count - 1)
$targetPos = $lastContentParagraph.Range.End - 1

$placeholderIns = $d.Range($targetPos, $targetPos)
$placeholderIns.InsertAfter("X")

$bookmarkRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholder = $d.Range($targetPos, $targetPos + 1)
$placeholder.Delete()

$trailingParagraph = $d.Paragraphs($d.Paragraphs.Count)
$trailingParagraph.Range.Delete()
